$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.034.42"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "1.832.87"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6710"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.22%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07416"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2952"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07650"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("D12").Value = "1.838.19"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.005"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6725"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.137"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.46%  "
$ws.Range("D17").Value = "29.025.30"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008230"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "2.071.07"
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.21%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.316"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.83%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9995"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1429"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.56%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.667"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.82%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.500"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.229"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.110"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.195"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05368"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.10%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7495"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.71%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.852"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.122"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.65%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.683"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.291.61"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01808"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.43%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.709"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9222"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.83%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.036"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.74%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9988"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("B45").Value = "XinFinNetwork"
$ws.Range("C45").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.08225"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +25.55%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000128"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.18%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.972.63"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5175"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.75%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "63.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.348"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.751"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.57%  "
